$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder info
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 10.12.2023"

# Row 6
$ws.Range("B6").Value = "12.12."
$ws.Range("C6").Value = "13.12."
$ws.Range("D6").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E6").Value = "57,95-"

# Row 7
$ws.Range("B7").Value = "16.12."
$ws.Range("C7").Value = "17.12."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-67927348"
$ws.Range("E7").Value = "52,73-"

# Row 8
$ws.Range("B8").Value = "19.12."
$ws.Range("C8").Value = "20.12."
$ws.Range("D8").Value = "MCDONALDS Gera"
$ws.Range("E8").Value = "31,83-"

# Row 9
$ws.Range("B9").Value = "20.12."
$ws.Range("C9").Value = "21.12."
$ws.Range("D9").Value = "PAYPAL VZJBVK"
$ws.Range("E9").Value = "89,42-"

# Row 10
$ws.Range("B10").Value = "24.12."
$ws.Range("C10").Value = "25.12."
$ws.Range("D10").Value = "PAYPAL SQDPLP"
$ws.Range("E10").Value = "51,64-"

# Row 11 - cleared out (transaction removed)
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# Closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 27.12.2023"
$ws.Range("E12").Value = "283,57-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 01.01.2024"
